$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume }
$updates = @{
    2  = @{ D = "306.83";     E = "-6.45%" }
    3  = @{ D = "39.76";      E = "-10.34%" }
    4  = @{ D = "5.062";      E = "-6.34%" }
    5  = @{ D = "0.07772";    E = "-7.13%" }
    6  = @{ D = "4.323";      E = "-2.30%" }
    7  = @{ D = "1.646";      E = "-15.00%" }
    8  = @{ D = "0.9181";     E = "-5.63%" }
    9  = @{ D = "0.09686";    E = "-14.84%" }
    10 = @{ D = "0.1738";     E = "-8.60%" }
    11 = @{ D = "0.08972";    E = "-7.15%" }
    12 = @{ E = "-4.70%" }
    13 = @{ D = "7.041";      E = "-15.29%" }
    14 = @{ D = "0.1058";     E = "-0.30%" }
    15 = @{ D = "0.001275";   E = "-1.34%" }
    16 = @{ D = "0.005661";   E = "-5.11%" }
    17 = @{ D = "3.367";      E = "0.07%" }
    18 = @{ D = "2.580";      E = "2.78%" }
    19 = @{ D = "0.3367";     E = "0.29%" }
    20 = @{ D = "0.1367";     E = "-1.54%" }
    21 = @{ D = "0.2665";     E = "0.42%" }
    22 = @{ D = "0.04145";    E = "-0.72%" }
    23 = @{ D = "0.001205";   E = "-2.57%" }
    24 = @{ D = "0.004086";   E = "-8.17%" }
    25 = @{ D = "0.0001227";  E = "-5.55%" }
    26 = @{ D = "0.0002997";  E = "0.61%" }
    38 = @{ D = "0.02378";    E = "-12.30%" }
    39 = @{ D = "0.05156";    E = "-8.29%" }
    40 = @{ D = "0.007980";   E = "1.62%" }
    41 = @{ E = "-6.17%" }
    42 = @{ D = "0.007563";   E = "3.74%" }
    43 = @{ D = "0.002020";   E = "-0.95%" }
    44 = @{ D = "0.008071";   E = "-7.35%" }
    45 = @{ D = "0.3330";     E = "-4.98%" }
    46 = @{ D = "0.00006744"; E = "-2.31%" }
    47 = @{ E = "0.64%" }
    48 = @{ D = "0.003428";   E = "-1.59%" }
    49 = @{ D = "0.004125";   E = "16.86%" }
    50 = @{ D = "0.00002113"; E = "0.64%" }
    51 = @{ D = "0.0002012";  E = "0.64%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["D"]
    }
    if ($cols.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["E"]
    }
}
